# edit.ps1 - applies the tracked changes described by the diff:
#   1. "Introduction" title: font size 20pt -> 12pt (w:sz 40 -> 24), bold kept
#   2. "Purpose of the system" heading: remove bold
#   3. "Scope of the system" heading: remove bold
#   4. "Objectives and success criteria of the project" heading: remove bold
#   5. "Perform authentication on u|sers..." : merge the two split runs (and
#      the _GoBack bookmark that sat between them) back into a single run
#   6. Re-insert the _GoBack bookmark, now collapsed, right before the
#      "References" heading run

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $exactText) {
    # Locate a paragraph whose trimmed text matches $exactText exactly -
    # more robust than hard-coded paragraph indices.
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.TrimEnd("`r`a`n") -eq $exactText) {
            return $p
        }
    }
    return $null
}

# --- 1. Shrink the "Introduction" title from 20pt to 12pt -----------------
# Selecting the whole paragraph Range (which includes the paragraph mark)
# so both the run and the paragraph-mark rPr pick up the new size, same as
# Word does when you select the whole line and change the font size.
$titlePara = Get-ParagraphByText $d "Introduction"
$titlePara.Range.Font.Size = 12

# --- 2-4. Un-bold the three section headings -------------------------------
(Get-ParagraphByText $d "Purpose of the system").Range.Font.Bold = $false
(Get-ParagraphByText $d "Scope of the system").Range.Font.Bold = $false
(Get-ParagraphByText $d "Objectives and success criteria of the project").Range.Font.Bold = $false

# --- 5. Merge the "Perform authentication on u" / "sers before..." runs ---
# The original paragraph is split into two runs around a _GoBack bookmark:
#   "Perform authentication on u" + <bookmark> + "sers before allowing any access to the system"
# Replacing the whole paragraph text with itself collapses it back into a
# single run and removes the bookmark that was sitting inside the range.
$authText = "Perform authentication on users before allowing any access to the system"
$authParRange = (Get-ParagraphByText $d $authText).Range
$authParRange.Find.Execute($authText, $true, $false, $false, $false, $false, $true, 1, $false, $authText, 2) | Out-Null

# --- 6. Re-add the _GoBack bookmark just before the "References" heading --
# Collapse a zero-length range to the start of the "References" run, then
# add the (hidden) _GoBack bookmark there - this is what Word leaves behind
# marking the last edit position in a session.
$refRange = $d.Content
$refRange.Find.Execute("References", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $refRange.Start
$collapsed = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $collapsed) | Out-Null

Write-Output "edits applied"
